$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 300 (existing rows 300:374 shift down to 302:376).
$ws.Rows("300:301").Insert()

# New row 300 data
$ws.Range("A300").Value = 4
$ws.Range("B300").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C300").Value = "Los Lagos"
$ws.Range("D300").Value = 44627
$ws.Range("E300").Value = 10
$ws.Range("F300").Value = "Fruta"
$ws.Range("G300").Value = 100106
$ws.Range("H300").Value = "Oleaginosos"
$ws.Range("I300").Value = 100106002
$ws.Range("J300").Value = "Palta"
$ws.Range("K300").Value = "Hass"
$ws.Range("L300").Value = "Primera"
$ws.Range("M300").Value = 200
$ws.Range("N300").Value = 4300
$ws.Range("O300").Value = 4300
$ws.Range("P300").Value = 4300
$ws.Range("Q300").Value = "$/kilo (en caja de 17 kilos)"
$ws.Range("R300").Value = "Provincia de Quillota"
$ws.Range("S300").Value = 4300
$ws.Range("T300").Value = 1

# New row 301 data
$ws.Range("A301").Value = 4
$ws.Range("B301").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C301").Value = "Los Lagos"
$ws.Range("D301").Value = 44627
$ws.Range("E301").Value = 10
$ws.Range("F301").Value = "Fruta"
$ws.Range("G301").Value = 100106
$ws.Range("H301").Value = "Oleaginosos"
$ws.Range("I301").Value = 100106002
$ws.Range("J301").Value = "Palta"
$ws.Range("K301").Value = "Hass"
$ws.Range("L301").Value = "Segunda"
$ws.Range("M301").Value = 100
$ws.Range("N301").Value = 3900
$ws.Range("O301").Value = 3900
$ws.Range("P301").Value = 3900
$ws.Range("Q301").Value = "$/kilo (en caja de 17 kilos)"
$ws.Range("R301").Value = "Provincia de Quillota"
$ws.Range("S301").Value = 3900
$ws.Range("T301").Value = 1
